# Append the new match row (row 3) for Arshdeep Singh's stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = $ws.Range("A3:K3")

$ws.Range("A3").Value = " Dubai (DSC)"
$ws.Range("B3").Value = " October 08 2020"
$ws.Range("C3").Value = "Sunrisers won by 69 runs"
$ws.Range("D3").Value = "Kings XI Punjab"
$ws.Range("E3").Value = "Sunrisers Hyderabad"
$ws.Range("F3").Value = "Arshdeep Singh" + [char]0x00A0

# These look numeric but must stay text (like row 2) - a leading apostrophe
# forces text storage; ClearFormats() afterwards drops the quote-prefix
# style Excel applies so the cells end up with the sheet's default style,
# matching the plain (unstyled) text cells already on the sheet.
$ws.Range("G3").Value = "'0"
$ws.Range("H3").Value = "'3"
$ws.Range("I3").Value = "'0"
$ws.Range("J3").Value = "'0"
$ws.Range("K3").Value = "'0.00"

$newRow.ClearFormats()
